$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$cell = $ws.Range("C1")
$cell.NumberFormat = "mm-dd-yy"
$cell.Value = (Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0).Date
